$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
Write-Host "F1 MergeCells: $($ws2.Range('F1').MergeCells)"
Write-Host "F1 MergeArea: $($ws2.Range('F1').MergeArea.Address())"
